$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Languages: insert "Swift 4 (proficient), " between
#    "C++ (proficient), " and "C# (proficient..."
# ------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("C# (profi", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(1)
$r.InsertBefore("Swift 4 (proficient), ")

# ------------------------------------------------------------------
# 2) Languages: "Swift 4 (basics), Dart (basics)" -> "Dart (basic)"
# ------------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute("Swift 4 (basics), Dart (basics)", $true, $false, $false, $false, $false, $true, 1, $false, "Dart (basic)", 1)

# ------------------------------------------------------------------
# 3) Move the _GoBack bookmark from after "commerce applications that
#    run" to inside "Moshi" (between "Mo" and "shi")
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$r3 = $d.Content
$r3.Find.Execute("Moshi", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$mid = $d.Range($r3.Start + 2, $r3.Start + 2)
$d.Bookmarks.Add("_GoBack", $mid)

# ------------------------------------------------------------------
# 4) Computer Village Hub: " (2020)" -> " (February 2020 - Present)"
# ------------------------------------------------------------------
$r4 = $d.Content
$r4.Find.Execute("www.cvhub4africa.com", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r4.Collapse(0)
$r4.Find.Execute(" (2020)", $true, $false, $false, $false, $false, $true, 1, $false, " (February 2020 - Present)", 1)

# ------------------------------------------------------------------
# 5) Law Pavilion hyperlink: " (2018)" -> " (June 2018 – November 2018)"
# ------------------------------------------------------------------
$r5 = $d.Content
$r5.Find.Execute("Law Pavilion, Grace InfoTech Limited", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r5.Collapse(0)
$r5.Find.Execute(" (2018)", $true, $false, $false, $false, $false, $true, 1, $false, " (June 2018 – November 2018)", 1)
